# "No MyKad mesti unik" - mark the "MyKad no mesti unik" bullet item as
# struck-through (done), leaving everything else in the document untouched.

$d = $word.ActiveDocument

# Locate the paragraph that reads "MyKad no mesti unik;" and apply
# strike-through formatting to the whole paragraph (text + paragraph mark),
# matching how Word records a manual strike-through edit in the OOXML
# (a <w:strike/> element added to every run's rPr as well as the
# paragraph-mark rPr).
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $text = $para.Range.Text
    if ($text -like "MyKad no mesti unik*") {
        $para.Range.Font.StrikeThrough = 1
        break
    }
}
